$d = $word.ActiveDocument

# --- 1) Mark the run that holds the final "ER" image as noProof -----------
$lastParaIndex = $d.Paragraphs.Count
$erImagePara = $d.Paragraphs.Item($lastParaIndex)
$erImagePara.Range.NoProofing = 1

# --- 2) Add the new content after the ER-image paragraph ------------------
# First create a fresh paragraph right after the image; it inherits the
# Arial/24pt run formatting that is already in effect at that point.
$erImagePara.Range.InsertParagraphAfter()

$newParaIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newParaIndex)

$bodyXml = "<w:p>" +
    "<w:r><w:rPr><w:rFonts w:ascii=""Arial"" w:hAnsi=""Arial"" w:cs=""Arial""/><w:sz w:val=""24""/><w:szCs w:val=""24""/></w:rPr><w:br w:type=""page""/></w:r>" +
    "<w:r><w:rPr><w:rFonts w:ascii=""Arial"" w:hAnsi=""Arial"" w:cs=""Arial""/><w:sz w:val=""24""/><w:szCs w:val=""24""/></w:rPr><w:lastRenderedPageBreak/><w:t>C</w:t></w:r>" +
    "<w:r><w:t>onsultas en forma SQL y su resultado.</w:t></w:r>" +
    "</w:p>"

$packageXml = "<pkg:package xmlns:pkg=""http://schemas.microsoft.com/office/2006/xmlPackage"">" +
    "<pkg:part pkg:name=""/word/document.xml"" pkg:contentType=""application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml""><pkg:xmlData>" +
    "<w:document xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"">" +
    "<w:body>" + $bodyXml + "</w:body></w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

# InsertXML replaces the content of the (still empty) paragraph we just
# created. Word keeps the paragraph's own mark as a trailing, run-less
# paragraph, which is exactly the blank "Arial 24pt" paragraph that closes
# out the document.
$newPara.Range.InsertXML($packageXml)

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
